# Apply the "commit" update:
#  1. Column C ("Förändrad") date serial is bumped from 45184 to 45186 for every data row.
#  2. The HYPERLINK() formulas in columns S, T, V, W, X, Y (rows 2-9, the only rows that
#     contain them) gain a second argument - the friendly link text - which is simply the
#     case identifier found in column A of the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq 45184) {
        $cCell.Value2 = 45186
    }

    $label = $ws.Cells.Item($r, 1).Value2
    if ($label) {
        foreach ($col in 19, 20, 22, 23, 24, 25) {
            $cell = $ws.Cells.Item($r, $col)
            $formula = $cell.Formula
            if ($formula -like 'HYPERLINK(*' -or $formula -like '=HYPERLINK(*') {
                if ($formula -notlike '*,*') {
                    $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $label + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
